# The "BESPOKE GRAMMAR SCHOOL OF ENGLISH" promotion has been renamed to
# "ENGLISH" -- update the department column (C) for both data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ENGLISH"
$ws.Range("C3").Value = "ENGLISH"

# Leave the selection where the author's last save left it.
$ws.Range("C3").Select()
